# Auto-generated script applying cell-value updates per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Range('D2').Value = '34.692.27'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '1.793.28'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.552'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.42'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.283'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0712'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.36%  '
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').Value = '2.051.05'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.797.87'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D16').Value = '34.718.81'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('E17').Value = '  +1.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '254.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = '0.0₃0807'
$ws.Range('E20').Value = '  +8.21%  '
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.75'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0529'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.90%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').Value = '1.439.50'
$ws.Range('E35').Value = '  -3.30%  '
$ws.Range('E36').Value = '  +2.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.634'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '84.73'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.927'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.50%  '
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('E46').Value = '  -4.71%  '
$ws.Range('D47').Value = '1.948.64'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.77%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.97'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('E51').Value = '  +7.91%  '
